# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-10-03 (serial 45202) to 2023-10-04 (serial 45203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 151; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}
